$d = $word.ActiveDocument

# 1) "To add new timer user must specify " + "timeout" + ", timer mode, ..."
#    were split across three runs with identical content/formatting; collapse
#    them into a single run carrying the same, unchanged sentence (pure
#    run-merge / typo cleanup, no wording change).
$d.Content.Find.Execute(
    "To add new timer user must specify timeout, timer mode, timeout resolution, pointer to callback function and pointer to data that should be passed to callback function.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To add new timer user must specify timeout, timer mode, timeout resolution, pointer to callback function and pointer to data that should be passed to callback function.",
    2) | Out-Null

# 2) Remove the stray manual line break between "...milliseconds. " and
#    "Timeout must be greater..." so the two sentences share one line/run.
#    "^l" is the Find special-character code for a manual line break
#    (w:br) and is recognised even with MatchWildcards off.
$d.Content.Find.Execute(
    "Timeout resolution is specified as seconds or milliseconds. ^lTimeout must be greater than zero for specified timer mode.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Timeout resolution is specified as seconds or milliseconds. Timeout must be greater than zero for specified timer mode.",
    2) | Out-Null

# 3) "...set on system power cycle and when HW Timer started/stopped."
#    becomes "...set on system power cycle and with HW Timer start/stop."
#    Insert the new wording right before the _GoBack bookmark (so the
#    bookmark keeps sitting right after "and"), then trim the old tail
#    down to a bare period.
$bm = $d.Bookmarks("_GoBack")
$ins = $bm.Range
$ins.InsertAfter(" w")
$ins.InsertAfter("ith")
$ins.InsertAfter(" HW Timer start/stop")

$d.Content.Find.Execute(
    " when HW Timer started/stopped.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ".",
    2) | Out-Null
